$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force column A to be entered as text so the dd-mm-yyyy strings are not
# auto-converted into Excel date serial numbers.
$ws.Range("A148:A149").NumberFormat = "@"

# Row 148
$ws.Range("A148").Value = "03-08-2021"
$ws.Range("B148").Value = 4050
$ws.Range("C148").Value = 1031
$ws.Range("D148").Value = 734
$ws.Range("E148").Value = 652
$ws.Range("F148").Value = 698
$ws.Range("G148").Value = 935

# Row 149
$ws.Range("A149").Value = "04-08-2021"
$ws.Range("B149").Value = 4958
$ws.Range("C149").Value = 1093
$ws.Range("D149").Value = 615
$ws.Range("E149").Value = 818
$ws.Range("F149").Value = 1170
$ws.Range("G149").Value = 1261

# Restore the default cell style so the new rows match the formatting of
# the rest of the data (no explicit style index), same as before entering
# the text-forced number format above.
$ws.Range("A148:A149").Style = "Normal"
